$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 139.84616
$ws.Range("I5").Value = 129.3
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 129.3
$ws.Range("L5").Value = 175
$ws.Range("M5").Value = -14.30000000000001
$ws.Range("N5").Value = -405
# Row 12
$ws.Range("H12").Value = 905.5
$ws.Range("I12").Value = 795.5
$ws.Range("J12").Value = 1125.5
$ws.Range("K12").Value = 795.5
$ws.Range("L12").Value = 1125.5
$ws.Range("M12").Value = -625.5
$ws.Range("N12").Value = -1465.5
# Row 55
$ws.Range("H55").Value = 836.3889
$ws.Range("I55").Value = 713.1667
$ws.Range("J55").Value = 898
$ws.Range("K55").Value = 713.1667
$ws.Range("L55").Value = 898
$ws.Range("M55").Value = -499.1667
$ws.Range("N55").Value = -1326
# Row 88
$ws.Range("H88").Value = 1999.4
$ws.Range("J88").Value = 1999
$ws.Range("L88").Value = 1999
$ws.Range("N88").Value = -2811
# Row 91
$ws.Range("H91").Value = 1999.4
$ws.Range("J91").Value = 1999
$ws.Range("L91").Value = 1999
$ws.Range("N91").Value = -4807
# Row 101
$ws.Range("H101").Value = 419.75
$ws.Range("J101").Value = 700
$ws.Range("L101").Value = 2100
$ws.Range("N101").Value = -5344

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 4996.143
$ws.Range("I63").Value = 2618.75
$ws.Range("J63").Value = 8166
$ws.Range("K63").Value = 2618.75
$ws.Range("L63").Value = 8166
$ws.Range("M63").Value = -1932.75
$ws.Range("N63").Value = -9538
# Row 66
$ws.Range("H66").Value = 4996.143
$ws.Range("I66").Value = 2618.75
$ws.Range("J66").Value = 8166
$ws.Range("K66").Value = 13093.75
$ws.Range("L66").Value = 40830
$ws.Range("M66").Value = -9661.75
$ws.Range("N66").Value = -47694
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 97
$ws.Range("H97").Value = 732.2778
$ws.Range("I97").Value = 772.75
$ws.Range("K97").Value = 772.75
$ws.Range("M97").Value = -276.75
# Row 122
$ws.Range("H122").Value = 1561.625
$ws.Range("I122").Value = 1520.7858
$ws.Range("K122").Value = 4562.357400000001
$ws.Range("M122").Value = -2112.357400000001
# Row 132
$ws.Range("H132").Value = 7663.5
$ws.Range("I132").Value = 8836.200000000001
$ws.Range("K132").Value = 26508.6
$ws.Range("M132").Value = -23978.6

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 636.7143
$ws.Range("I94").Value = 576.3333
$ws.Range("J94").Value = 999
$ws.Range("K94").Value = 576.3333
$ws.Range("L94").Value = 999
$ws.Range("M94").Value = -125.3333
$ws.Range("N94").Value = -1901
# Row 105
$ws.Range("H105").Value = 2010
$ws.Range("I105").Value = 2010
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2010
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -263
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 2022
$ws.Range("I122").Value = 1559.3636
$ws.Range("J122").Value = 4566.5
$ws.Range("K122").Value = 4678.0908
$ws.Range("L122").Value = 13699.5
$ws.Range("M122").Value = -2228.0908
$ws.Range("N122").Value = -18599.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1126.5714
$ws.Range("J5").Value = 1549.5
$ws.Range("L5").Value = 4648.5
$ws.Range("N5").Value = -4872.5
# Row 11
$ws.Range("H11").Value = 1219.8
$ws.Range("I11").Value = 366.66666
$ws.Range("K11").Value = 1099.99998
$ws.Range("M11").Value = -959.9999800000001
# Row 82
$ws.Range("H82").Value = 715
$ws.Range("J82").Value = 715
$ws.Range("L82").Value = 2145
$ws.Range("N82").Value = -2957
# Row 85
$ws.Range("H85").Value = 715
$ws.Range("J85").Value = 715
$ws.Range("L85").Value = 2145
$ws.Range("N85").Value = -4953
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
# Row 103
$ws.Range("H103").Value = 549
$ws.Range("I103").Value = 499.5
$ws.Range("J103").Value = 598.5
$ws.Range("K103").Value = 1498.5
$ws.Range("L103").Value = 1795.5
$ws.Range("M103").Value = -619.5
$ws.Range("N103").Value = -3553.5
# Row 122
$ws.Range("H122").Value = 903
$ws.Range("J122").Value = 997.2
$ws.Range("L122").Value = 8974.800000000001
$ws.Range("N122").Value = -13874.8
# Row 131
$ws.Range("H131").Value = 1333
$ws.Range("J131").Value = 1333
$ws.Range("L131").Value = 3999
$ws.Range("N131").Value = -14079
# Row 135
$ws.Range("H135").Value = 1126.5714
$ws.Range("J135").Value = 1549.5
$ws.Range("L135").Value = 13945.5
$ws.Range("N135").Value = -19015.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 239.05882
$ws.Range("I2").Value = 184.11111
$ws.Range("K2").Value = 184.11111
$ws.Range("M2").Value = -71.11111
# Row 3
$ws.Range("H3").Value = 18105518
$ws.Range("I3").Value = 26825396
$ws.Range("J3").Value = 12001602
$ws.Range("K3").Value = 26825396
$ws.Range("L3").Value = 12001602
$ws.Range("M3").Value = -26825280
$ws.Range("N3").Value = -12001834

$ws = $wb.Worksheets.Item("LTW")
# Row 70
$ws.Range("H70").Value = 20949.5
$ws.Range("J70").Value = 20949.5
$ws.Range("L70").Value = 20949.5
$ws.Range("N70").Value = -21489.5
# Row 73
$ws.Range("H73").Value = 20949.5
$ws.Range("J73").Value = 20949.5
$ws.Range("L73").Value = 20949.5
$ws.Range("N73").Value = -22821.5
# Row 93
$ws.Range("H93").Value = 1588.8572
$ws.Range("I93").Value = 2032.25
$ws.Range("J93").Value = 997.6667
$ws.Range("K93").Value = 2032.25
$ws.Range("L93").Value = 997.6667
$ws.Range("M93").Value = -784.25
$ws.Range("N93").Value = -3493.6667
# Row 122
$ws.Range("H122").Value = 3966.1667
$ws.Range("I122").Value = 3966.1667
$ws.Range("K122").Value = 11898.5001
$ws.Range("M122").Value = -9448.500100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 34545.453
$ws.Range("I2").Value = 34545.453
$ws.Range("K2").Value = 34545.453
$ws.Range("M2").Value = -34433.453
# Row 122
$ws.Range("H122").Value = 4950.5
$ws.Range("I122").Value = 3267.3333
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 9801.999899999999
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -7351.999899999999
$ws.Range("N122").Value = -34900
# Row 136
$ws.Range("H136").Value = 2449.2144
$ws.Range("I136").Value = 1521.2727
$ws.Range("J136").Value = 5851.6665
$ws.Range("K136").Value = 4563.8181
$ws.Range("L136").Value = 17554.9995
$ws.Range("M136").Value = -2013.8181
$ws.Range("N136").Value = -22654.9995
